$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised data values (columns B, C, D) across existing rows
$ws.Range("C68").Value = 83.7
$ws.Range("C71").Value = 82.40000000000001
$ws.Range("C74").Value = 81
$ws.Range("C75").Value = 80.8
$ws.Range("C85").Value = 82.7
$ws.Range("D98").Value = 1.6
$ws.Range("D100").Value = 1.3
$ws.Range("D108").Value = 1.1
$ws.Range("C112").Value = 95.2
$ws.Range("C115").Value = 96.2
$ws.Range("D123").Value = 1.3
$ws.Range("D136").Value = 0.2
$ws.Range("D137").Value = 0.5
$ws.Range("D152").Value = 0.9
$ws.Range("C156").Value = 104.7
$ws.Range("D159").Value = 0.6
$ws.Range("C160").Value = 105.8
$ws.Range("D165").Value = 0.3
$ws.Range("D166").Value = 0.3
$ws.Range("D168").Value = 0.4
$ws.Range("D170").Value = 0.1
$ws.Range("C173").Value = 105.8
$ws.Range("D174").Value = 0
$ws.Range("D175").Value = 1.1
$ws.Range("C176").Value = 107.7
$ws.Range("C177").Value = 108
$ws.Range("D177").Value = 1.8
$ws.Range("D178").Value = 1.7
$ws.Range("C182").Value = 110.1
$ws.Range("C183").Value = 110.4
$ws.Range("C185").Value = 111.4
$ws.Range("C188").Value = 111.2
$ws.Range("D188").Value = 0.4
$ws.Range("D189").Value = -0.3
$ws.Range("C190").Value = 111.2
$ws.Range("D190").Value = -0.4
$ws.Range("D191").Value = -0.2
$ws.Range("D193").Value = 1.2
$ws.Range("D194").Value = 1.2
$ws.Range("D195").Value = -0.2
$ws.Range("D196").Value = -0.2
$ws.Range("C197").Value = 113.2
$ws.Range("D197").Value = -0.1
$ws.Range("C199").Value = 112.9
$ws.Range("D199").Value = 1.2
$ws.Range("C200").Value = 113.1
$ws.Range("D201").Value = 0.4
$ws.Range("D202").Value = 0.6
$ws.Range("D203").Value = -0.9
$ws.Range("C205").Value = 112
$ws.Range("C206").Value = 113.8
$ws.Range("D206").Value = -0.6
$ws.Range("C207").Value = 114
$ws.Range("D207").Value = 2.7
$ws.Range("C208").Value = 107.7
$ws.Range("C209").Value = 99.2
$ws.Range("C211").Value = 96.90000000000001
$ws.Range("C212").Value = 99
$ws.Range("D213").Value = -2
$ws.Range("C214").Value = 107
$ws.Range("D214").Value = 5.1
$ws.Range("C215").Value = 107.4
$ws.Range("D215").Value = 8
$ws.Range("D217").Value = 6.6
$ws.Range("B218").Value = 106.7
$ws.Range("B219").Value = 104.2
$ws.Range("C219").Value = 114.1
$ws.Range("D219").Value = 4.7
$ws.Range("B220").Value = 120.5
$ws.Range("C220").Value = 113
$ws.Range("D220").Value = 3.5
$ws.Range("B221").Value = 112.8
$ws.Range("C221").Value = 111.5
$ws.Range("D221").Value = 1.4
$ws.Range("B222").Value = 116.6
$ws.Range("C222").Value = 115.2
$ws.Range("D222").Value = 0.2
$ws.Range("B223").Value = 115.9
$ws.Range("C223").Value = 117.2
$ws.Range("D223").Value = 1.2

# Append new row 224 for the 01-07-2021 period
$ws.Range("A224").NumberFormat = "@"
$ws.Range("A224").Value = "01-07-2021"
$ws.Range("A224").Style = "Normal"
$ws.Range("B224").Value = 113.5
$ws.Range("C224").Value = 118.8
$ws.Range("D224").Value = 3.7

